$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text format so numeric-looking strings
# (e.g. "319.35", "43.095.82", "  -4.85%  ") are preserved verbatim as text
# instead of being auto-converted into numbers by Excel.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '43.095.82'
$ws.Range('E2').Value = '  -4.85%  '
$ws.Range('D3').Value = '2.228.17'
$ws.Range('E3').Value = '  -5.92%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '319.35'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').Value = '99.55'
$ws.Range('E6').Value = '  -8.28%  '
$ws.Range('D7').Value = '0.583'
$ws.Range('E7').Value = '  -8.29%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '0.560'
$ws.Range('E9').Value = '  -8.60%  '
$ws.Range('D10').Value = '36.91'
$ws.Range('E10').Value = '  -9.95%  '
$ws.Range('D11').Value = '54.14'
$ws.Range('E11').Value = '  -3.45%  '
$ws.Range('D12').Value = '0.0830'
$ws.Range('E12').Value = '  -9.58%  '
$ws.Range('D13').Value = '7.66'
$ws.Range('E13').Value = '  -9.73%  '
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.569.95'
$ws.Range('E15').Value = '  -5.94%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '0.863'
$ws.Range('E16').Value = '  -12.06%  '
$ws.Range('D17').Value = '14.38'
$ws.Range('E17').Value = '  -6.62%  '
$ws.Range('D18').Value = '2.238.17'
$ws.Range('E18').Value = '  -5.80%  '
$ws.Range('D19').Value = '43.028.37'
$ws.Range('E19').Value = '  -5.04%  '
$ws.Range('D20').Value = '14.55'
$ws.Range('E20').Value = '  -4.66%  '
$ws.Range('D21').Value = '0.0₃0965'
$ws.Range('E21').Value = '  -9.17%  '
$ws.Range('D22').Value = '6.53'
$ws.Range('E22').Value = '  -10.24%  '
$ws.Range('D23').Value = '65.11'
$ws.Range('E23').Value = '  -11.23%  '
$ws.Range('D24').Value = '3.16'
$ws.Range('E24').Value = '  -13.70%  '
$ws.Range('D25').Value = '237.47'
$ws.Range('E25').Value = '  -10.74%  '
$ws.Range('D26').Value = '2.16'
$ws.Range('E26').Value = '  -7.82%  '
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('D28').Value = '4.03'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('E29').Value = '  -2.05%  '
$ws.Range('D30').Value = '9.98'
$ws.Range('E30').Value = '  -10.63%  '
$ws.Range('D31').Value = '6.34'
$ws.Range('E31').Value = '  -15.30%  '
$ws.Range('D32').Value = '35.78'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('D33').Value = '20.33'
$ws.Range('E33').Value = '  -9.29%  '
$ws.Range('D34').Value = '0.0872'
$ws.Range('E34').Value = '  -7.85%  '
$ws.Range('D35').Value = '153.48'
$ws.Range('E35').Value = '  -9.22%  '
$ws.Range('D36').Value = '2.67'
$ws.Range('E36').Value = '  -5.07%  '
$ws.Range('D37').Value = '3.16'
$ws.Range('E37').Value = '  +6.16%  '
$ws.Range('D38').Value = '1.93'
$ws.Range('E38').Value = '  +1.40%  '
$ws.Range('E39').Value = '  -7.77%  '
$ws.Range('D40').Value = '4.43'
$ws.Range('E40').Value = '  -5.77%  '
$ws.Range('E41').Value = '  -11.13%  '
$ws.Range('D42').Value = '3.65'
$ws.Range('E42').Value = '  -8.08%  '
$ws.Range('D43').Value = '0.0323'
$ws.Range('E43').Value = '  -8.73%  '
$ws.Range('D44').Value = '13.49'
$ws.Range('E44').Value = '  +4.77%  '
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').Value = '1.742.17'
$ws.Range('E46').Value = '  -7.14%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.204'
$ws.Range('E47').Value = '  -10.00%  '
$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D48').Value = '85.13'
$ws.Range('E48').Value = '  -13.63%  '
$ws.Range('D49').Value = '5.30'
$ws.Range('E49').Value = '  -10.99%  '
$ws.Range('D50').Value = '75.20'
$ws.Range('E50').Value = '  -9.93%  '
$ws.Range('D51').Value = '8.70'
$ws.Range('E51').Value = '  -5.56%  '

# Restore the default cell style so no residual formatting is introduced
$priceRange.Style = "Normal"

